$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.074.51"
$ws.Range("E2").Value = "  +0.64%  "

$ws.Range("D3").Value = "2.927.35"
$ws.Range("E3").Value = "  +0.82%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.17%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("E11").Value = "  -1.85%  "

$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("E14").Value = "  -0.34%  "

$ws.Range("D15").Value = "3.410.21"
$ws.Range("E15").Value = "  +0.72%  "

$ws.Range("D16").Value = "61.031.43"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").Value = "2.926.30"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("E21").Value = "  +1.90%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("E25").Value = "  +0.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.97%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("D34").Value = "0.0₃0855"
$ws.Range("E34").Value = "  +2.32%  "

$ws.Range("E35").Value = "  +0.62%  "

$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.68%  "

$ws.Range("E38").Value = "  +1.67%  "

$ws.Range("E39").Value = "  -1.57%  "

$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "376.06"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "2.724.72"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0347"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "130.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.58%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("E48").Value = "  -3.22%  "

$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("E51").Value = "  +3.32%  "
